$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BranchData")
$ws.Range("S2:S64").NumberFormat = "#,##0"
$ws.Range("T2:V64").NumberFormat = "@"
$ws.Range("T2").Value = "600"
$ws.Range("U2").Value = "0"
$ws.Range("V2").Value = "716,000"
$ws.Range("T3").Value = "813"
$ws.Range("U3").Value = "30,000"
$ws.Range("V3").Value = "1,216,000"
$ws.Range("T4").Value = "0"
$ws.Range("U4").Value = "0"
$ws.Range("V4").Value = "577,000"
$ws.Range("T5").Value = "1,740"
$ws.Range("U5").Value = "0"
$ws.Range("V5").Value = "551,000"
$ws.Range("T6").Value = "304"
$ws.Range("U6").Value = "0"
$ws.Range("V6").Value = "1,581,000"
$ws.Range("T7").Value = "0"
$ws.Range("U7").Value = "0"
$ws.Range("V7").Value = "288,000"
$ws.Range("T8").Value = "764"
$ws.Range("U8").Value = "32,000"
$ws.Range("V8").Value = "790,000"
$ws.Range("T9").Value = "0"
$ws.Range("U9").Value = "250,700"
$ws.Range("V9").Value = "0"
$ws.Range("T10").Value = "699"
$ws.Range("U10").Value = "64,000"
$ws.Range("V10").Value = "1,689,669"
$ws.Range("T12").Value = "677"
$ws.Range("U12").Value = "25,750"
$ws.Range("V12").Value = "618,000"
$ws.Range("T14").Value = "0"
$ws.Range("U14").Value = "0"
$ws.Range("V14").Value = "439,000"
$ws.Range("V15").Value = "357,000"
$ws.Range("T16").Value = "0"
$ws.Range("U16").Value = "0"
$ws.Range("V16").Value = "258,000"
$ws.Range("T18").Value = "0"
$ws.Range("U18").Value = "0"
$ws.Range("V18").Value = "187,000"
$ws.Range("T19").Value = "250"
$ws.Range("U19").Value = "0"
$ws.Range("V19").Value = "262,000"
$ws.Range("T20").Value = "299"
$ws.Range("U20").Value = "30,000"
$ws.Range("V20").Value = "643,000"
$ws.Range("T21").Value = "630"
$ws.Range("U21").Value = "46,000"
$ws.Range("V21").Value = "205,000"
$ws.Range("T22").Value = "92"
$ws.Range("U22").Value = "14,225"
$ws.Range("V22").Value = "4,460,181"
$ws.Range("T23").Value = "6,992"
$ws.Range("U23").Value = "30,500"
$ws.Range("V23").Value = "1,110,000"
$ws.Range("T25").Value = "0"
$ws.Range("U25").Value = "3,000"
$ws.Range("V25").Value = "0"
$ws.Range("T29").Value = "1,564"
$ws.Range("U29").Value = "20,700"
$ws.Range("V29").Value = "0"
$ws.Range("T30").Value = "1,071"
$ws.Range("U30").Value = "20,400"
$ws.Range("V30").Value = "1,467,000"
$ws.Range("T31").Value = "730"
$ws.Range("U31").Value = "0"
$ws.Range("V31").Value = "260,000"
$ws.Range("T34").Value = "0"
$ws.Range("U34").Value = "48,000"
$ws.Range("V34").Value = "0"
$ws.Range("T35").Value = "1,492"
$ws.Range("U35").Value = "153,000"
$ws.Range("V35").Value = "1,266,000"
$ws.Range("T36").Value = "575"
$ws.Range("U36").Value = "36,000"
$ws.Range("V36").Value = "1,615,000"
$ws.Range("T42").Value = "0"
$ws.Range("U42").Value = "18,000"
$ws.Range("V42").Value = "1,825,000"
$ws.Range("T44").Value = "1,098"
$ws.Range("U44").Value = "36,300"
$ws.Range("V44").Value = "1,478,000"
$ws.Range("T47").Value = "0"
$ws.Range("U47").Value = "0"
$ws.Range("V47").Value = "436,000"
$ws.Range("T48").Value = "0"
$ws.Range("U48").Value = "9,500"
$ws.Range("V48").Value = "0"
$ws.Range("T49").Value = "763"
$ws.Range("U49").Value = "40,100"
$ws.Range("V49").Value = "1,413,000"
$ws.Range("T50").Value = "673"
$ws.Range("U50").Value = "22,000"
$ws.Range("V50").Value = "801,000"
$ws.Range("T54").Value = "0"
$ws.Range("U54").Value = "0"
$ws.Range("V54").Value = "2,178,000"
$ws.Range("T55").Value = "52"
$ws.Range("U55").Value = "0"
$ws.Range("V55").Value = "279,000"
$ws.Range("T56").Value = "0"
$ws.Range("U56").Value = "24,000"
$ws.Range("V56").Value = "0"
$ws.Range("T57").Value = "950"
$ws.Range("U57").Value = "0"
$ws.Range("V57").Value = "608,000"
$ws.Range("T58").Value = "1,210"
$ws.Range("U58").Value = "22,500"
$ws.Range("V58").Value = "1,267,820"
$ws.Range("T60").Value = "1,735"
$ws.Range("U60").Value = "42,000"
$ws.Range("V60").Value = "3,262,000"
$ws.Range("T61").Value = "21"
$ws.Range("U61").Value = "44,000"
$ws.Range("V61").Value = "0"
$ws.Range("T62").Value = "1,530"
$ws.Range("U62").Value = "0"
$ws.Range("V62").Value = "1,113,000"
$ws.Range("T63").Value = "0"
$ws.Range("U63").Value = "0"
$ws.Range("V63").Value = "1,985,000"
$ws.Range("T64").Value = "0"
$ws.Range("U64").Value = "0"
$ws.Range("V64").Value = "469,000"
$ws.Range("A1").Select()
